# Hortaliza, Femacal de La Calera - Sandia
# Insert a new weekly record row right before the current row 236,
# shifting all subsequent rows down by one (dimension grows from R261 to R262).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 236 (pushes old 236.."n" down by one).
$ws.Rows(236).Insert()

# Populate the newly inserted row 236 with the new weekly record.
$ws.Range("A236").Value = 3
$ws.Range("B236").Value = "Femacal de La Calera"
$ws.Range("C236").Value = "Coquimbo"
$ws.Range("D236").Value = 44491
$ws.Range("E236").Value = 5
$ws.Range("F236").Value = 100112028
$ws.Range("G236").Value = "Sandia"
$ws.Range("H236").Value = "Sin especificar"
$ws.Range("I236").Value = "Primera"
$ws.Range("J236").Value = 180
$ws.Range("K236").Value = 800
$ws.Range("L236").Value = 800
$ws.Range("M236").Value = 800
$ws.Range("N236").Value = '$/kilo (volumen en unidades)'
$ws.Range("O236").Value = "Perú"
$ws.Range("P236").Value = 800
$ws.Range("Q236").Value = 1
$ws.Range("R236").Value = "Hortaliza"
